# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "71.134.65"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.812.86"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.23%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "706.36"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "171.71"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.95%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.811.05"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").Value = "  -0.61%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.73"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +7.16%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.460"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000253"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.16%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "35.90"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.19%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.455.88"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.16%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.837.07"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.53%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "71.092.61"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "17.45"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.14"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("E20").Value = "  -0.13%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "501.95"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.69"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").Value = "  +0.25%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "84.22"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("E25").Value = "  -3.55%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.963.48"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.22%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.37"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  +0.20%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.04"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("E31").Value = "  -3.58%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.26"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  -2.48%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "29.05"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("E35").Value = "  -4.92%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "9.13"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.778.71"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.89%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  -2.63%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.38"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("E41").Value = "  -3.06%  "
$ws.Range("E42").Value = "  -1.80%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +0.12%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "166.78"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.84%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.000314"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "49.09"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "423.01"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.30%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.62"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("E51").Value = "  -2.40%  "
